$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A->B, B->C)
$ws.Range("A1").EntireColumn.Insert()

# Fill the new column A with sequential row numbers 1-10
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

# Update the selection to match the post-edit state (cell A11)
$ws.Range("A11").Select()
